# Weekly update: a new Coco price observation for
# "Mercado Mayorista Lo Valledor de Santiago" is inserted as the new first
# data record (row 26), pushing every existing observation down by one row
# (old row 26 -> 27, ... old row 71 -> 72). Sheet dimension grows from
# A1:T71 to A1:T72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 26; rows 26-71 shift down to become rows 27-72.
$ws.Rows.Item(26).Insert()

# Seed the new row 26 by duplicating the (now shifted-down) row 27, which
# carries all the fixed descriptive columns (market, region, product,
# category, unit, origin, kg/unit, etc.) that stay identical for the new
# observation.
$ws.Range("A27:T27").Copy($ws.Range("A26"))

# Overwrite the columns that actually differ for this new observation:
# date, volume, min/max/weighted-avg price and $/kg.
$ws.Range("D26").Value2 = 44775
$ws.Range("M26").Value2 = 150
$ws.Range("N26").Value2 = 22000
$ws.Range("O26").Value2 = 22000
$ws.Range("P26").Value2 = 22000
$ws.Range("S26").Value2 = 1100
